# Add new resource rows (Coin / Crystal / Flower / Leaf / Money / Other) with
# their initial quantities on the "resource" sheet, and move the active
# selection to E15, matching the authored commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("resource")

$ws.Range("B5").Value = "Coin"
$ws.Range("E5").Value = 1000

$ws.Range("B6").Value = "Crystal"
$ws.Range("E6").Value = 100

$ws.Range("B7").Value = "Flower"
$ws.Range("E7").Value = 0

$ws.Range("B8").Value = "Leaf"
$ws.Range("E8").Value = 180

$ws.Range("B9").Value = "Money"

$ws.Range("B10").Value = "Other"

$ws.Range("E15").Select()
